# Update workbook per release-notes.md regeneration (ror-contact-description IG):
#  - Metadata sheet: Version / Status / Date / Contact values
#  - Elements sheet: swap the two "Mapping" columns (AK <-> AL), header + data + widths

$wb = $excel.ActiveWorkbook

# --- Metadata sheet -------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

$meta.Range("B3").Value  = "0.4.0-snapshot-1"                 # Version
$meta.Range("B6").Value  = "draft"                             # Status
$meta.Range("B8").Value  = "2024-05-23T12:16:26+00:00"         # Date
$meta.Range("B10").Value = "ANS (https://esante.gouv.fr)"      # Contact

# --- Elements sheet: swap Mapping columns (AK <-> AL) ----------------------
$elements = $wb.Worksheets.Item("Elements")

for ($r = 1; $r -le 6; $r++) {
    $akCell = $elements.Range("AK" + $r)
    $alCell = $elements.Range("AL" + $r)
    $akVal = $akCell.Value2
    $alVal = $alCell.Value2
    $akCell.Value = $alVal
    $alCell.Value = $akVal
}

# Swap the column widths to go with the swapped content (col 37 = AK, col 38 = AL)
$elements.Columns.Item(37).ColumnWidth = 73.33333333333334
$elements.Columns.Item(38).ColumnWidth = 24.166666666666664
